$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a full data row (columns A,B,D..K) in one shot. Column C
# ("N de compte") is handled separately below because it holds long,
# purely-numeric account numbers that must stay TEXT - if a cell's .Value
# is set to a literal digit string, the host coerces it straight to a
# Number (and a 24-25 digit account number loses its leading zero /
# precision in the process). Passing $null here leaves column C untouched
# for now so the row-write loop just skips it.
# ---------------------------------------------------------------------------
function Set-RowData {
    param($row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        if ($null -eq $values[$i]) { continue }
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Stage every distinct account-number string we will need as TEXT in
# far-away scratch cells *before* touching the real grid, while the
# original column C values are still pristine:
#   - 225400000805987601012173 already sits (as text) in C8 (unchanged by
#     this edit) - grab a copy.
#   - 145101211406073828000084 already sits (as text) in C3 - grab a copy
#     before row 3 is overwritten.
#   - 007400000313200019604463 is brand-new, so type it with a leading
#     apostrophe to force text.
# ---------------------------------------------------------------------------
$scratchAcct225 = $ws.Cells.Item(200, 1)
$scratchAcct145 = $ws.Cells.Item(200, 2)
$scratchAcct007 = $ws.Cells.Item(200, 3)

$ws.Cells.Item(8, 3).Copy()
$scratchAcct225.PasteSpecial(-4163)

$ws.Cells.Item(3, 3).Copy()
$scratchAcct145.PasteSpecial(-4163)

$scratchAcct007.Value = "'007400000313200019604463"

# ---------------------------------------------------------------------------
# Row-by-row target data for columns A,B,D..K (column C = $null, filled in
# further down from the scratch cells staged above).
# ---------------------------------------------------------------------------
Set-RowData 2  @("NOUBAIL MOHAMMED","IR801997",$null,"KHOURIBGA ZELLAKA","AWB","Direction régionale","444/444","mensuelle",22500,0,22500)
Set-RowData 3  @("JEMAA HORMI","B219321",$null,"KHOURIBGA","CA","Direction régionale","444/444","mensuelle",22500,0,22500)
Set-RowData 4  @("NOUBAIL MOHAMMED","IR801997",$null,"KHOURIBGA ZELLAKA","AWB","Direction régionale","444/444","mensuelle",3000,0,3000)
Set-RowData 5  @("JEMAA HORMI","B219321",$null,"KHOURIBGA","CA","Direction régionale","444/444","mensuelle",3000,0,3000)
Set-RowData 6  @("NOUBAIL MOHAMMED","IR801997",$null,"KHOURIBGA ZELLAKA","AWB","Direction régionale","444/444","mensuelle",1500,0,1500)
Set-RowData 7  @("JEMAA HORMI","B219321",$null,"KHOURIBGA","CA","Direction régionale","444/444","mensuelle",1500,0,1500)
Set-RowData 8  @("DOUNIA LAMKADDAM","BK646476",$null,"KHOURIBGA","CA","Direction régionale","000/CCCC/AV1","mensuelle",32000,0,32000)
Set-RowData 9  @("DOUNIA LAMKADDAM","BK646476",$null,"KHOURIBGA","CA","Direction régionale","000/CCCC/AV1","mensuelle",4000,0,4000)
Set-RowData 10 @("DOUNIA LAMKADDAM","BK646476",$null,"KHOURIBGA","CA","Direction régionale","000/CCCC/AV1","mensuelle",2000,0,2000)
Set-RowData 11 @("NADIA BADRANE","B171710",$null,"KHOURIBGA","CA","Direction régionale","555/RRR/AV10","mensuelle",2000,200,1800)
Set-RowData 12 @("NHILA BELGACEM","IB43905",$null,"MARRAKECH BENI MELLAL","BP","Direction régionale","555/RRR/AV10","mensuelle",2000,200,1800)
Set-RowData 13 @("NADIA BADRANE","B171710",$null,"KHOURIBGA","CA","Direction régionale","555/RRR/AV10","mensuelle",8500,850,7650)
Set-RowData 14 @("NHILA BELGACEM","IB43905",$null,"MARRAKECH BENI MELLAL","BP","Direction régionale","555/RRR/AV10","mensuelle",8500,850,7650)
Set-RowData 15 @("NADIA BADRANE","B171710",$null,"KHOURIBGA","CA","Direction régionale","555/RRR/AV10","mensuelle",2000,200,1800)
Set-RowData 16 @("NHILA BELGACEM","IB43905",$null,"MARRAKECH BENI MELLAL","BP","Direction régionale","555/RRR/AV10","mensuelle",2000,200,1800)
Set-RowData 17 @(" "," "," "," "," "," "," "," ",117000,2500,114500)

# ---------------------------------------------------------------------------
# Fan out the three text account numbers into column C without adding any
# per-cell number formatting (PasteSpecial xlPasteValues = -4163 keeps the
# TEXT value but not the source cell's style/quote-prefix formatting).
# ---------------------------------------------------------------------------
$scratchAcct007.Copy()
foreach ($r in @(2, 4, 6)) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4163)
}

$scratchAcct225.Copy()
foreach ($r in @(3, 5, 7, 8, 9, 10, 11, 13, 15)) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4163)
}

$scratchAcct145.Copy()
foreach ($r in @(12, 14, 16)) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4163)
}

$ws.Cells.Item(17, 3).Value = " "

# Clean up the scratch row used to stage the text account numbers.
$ws.Rows.Item(200).Delete()
